$wb = $excel.ActiveWorkbook

# --- "tech" sheet updates ---
$wsTech = $wb.Worksheets.Item("tech")

$wsTech.Range("D3").Value = 24.29
$wsTech.Range("F3").Value = 41.19

$wsTech.Range("D4").Value = 420
$wsTech.Range("F4").Value = 36.19

$wsTech.Range("D5").Value = 56
$wsTech.Range("F5").Value = 7.07

$wsTech.Range("D6").ClearContents() | Out-Null
$wsTech.Range("F6").Value = 25

$wsTech.Range("F7").Value = 7.45

$wsTech.Range("D11").Value = 344.7
$wsTech.Range("F11").Value = 8

# Page setup (paper size / orientation) similar to other sheets in the workbook
$wsTech.PageSetup.PaperSize = 9
$wsTech.PageSetup.Orientation = 1

# --- "co2" sheet updates ---
$wsCo2 = $wb.Worksheets.Item("co2")

$wsCo2.Range("D3").Value = 35
$wsCo2.Range("D4").Value = 275
$wsCo2.Range("D5").Value = 75
$wsCo2.Range("D6").Value = 505.13

# Update selection on co2 sheet to D11
$wsCo2.Range("D11").Select() | Out-Null

# Update selection on tech sheet to C24 and keep it as the active/visible tab
$wsTech.Activate() | Out-Null
$wsTech.Range("C24").Select() | Out-Null
